$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 3

$ws.Cells.Item($row, 1).Value = "JD_002"
$ws.Cells.Item($row, 2).Value = "Senior Engineer"
$ws.Cells.Item($row, 3).Value = "We are seeking a Software Engineer to build and maintain high-quality software solutions.`nWork with global teams to drive innovation and deliver scalable applications.`nJoin Akkodis and be part of a tech-driven, collaborative environment."
$ws.Cells.Item($row, 4).Value = 3
$ws.Cells.Item($row, 5).Value = 7
$ws.Cells.Item($row, 6).Value = "Remote"
$ws.Cells.Item($row, 7).Value = "Pune, Maharashtra, India"
